# Fig1.7_PV_module.pptx - "small details figures chapter 1"
#
# The only meaningful, reproducible content edit in the source diff is on
# slide 1: the caption reading "Back sheet (copolymer)" is split so the
# first word becomes "Backsheet" (no space) while the rest of the caption
# (" (copolymer)") is left as-is. (The diff's other hunks only touch the
# cached text of the auto-updating "datetimeFigureOut" date fields that
# live in the slide master / layouts; that is a side effect of PowerPoint
# recalculating those fields on save, not an addressable edit, so it is
# intentionally left untouched here.)

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# "TextBox 496" (shape #12 on slide 1) holds the caption text.
$shape = $s.Shapes.Item(12)
$tr = $shape.TextFrame.TextRange

# Replace() splits the run exactly at the match boundary, preserving the
# existing run formatting (Arial, 24pt) for both the replaced text and the
# untouched remainder.
$tr.Replace("Back sheet", "Backsheet") | Out-Null
